$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Edit")

# Column C (Skill Name) updates
$ws.Range("C3").Value = "ChatSkill11 <49011>"
$ws.Range("C4").Value = "AudioSkill-1 <49051>"
$ws.Range("C5").Value = "VideoSkill1 <49017>"
$ws.Range("C11").Value = "ChatSkill11 <49011>"
$ws.Range("C12").Value = "AudioSkill-1 <49051>"

# Column F (AgentList) updates
$ws.Range("F10").Value = "X Y - 1111,G A - 1213"
$ws.Range("F11").Value = "X Y - 1111,G A - 1213"
$ws.Range("F12").Value = "X Y - 1111,G A - 1213"
$ws.Range("F13").Value = "X Y - 1111,G A - 1213"
$ws.Range("F14").Value = "X Y - 1111,G A - 1213"
$ws.Range("F15").Value = "X Y - 1111,G A - 1213"
$ws.Range("F16").Value = "X Y - 1111,G A - 1213"
$ws.Range("F9").Value = "Aravinda ET - 602,Tester Web - 6189"

# Move the sheet's active selection cell to reflect the saved view state
$ws.Range("F10").Select()
